$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 35500
$ws.Range("J10").Value = 35500
$ws.Range("L10").Value = 35500
$ws.Range("N10").Value = -36086
$ws.Range("H100").Value = 22223658
$ws.Range("J100").Value = 1800
$ws.Range("L100").Value = 1800
$ws.Range("N100").Value = -2882
$ws.Range("H123").Value = 42743.332
$ws.Range("J123").Value = 42743.332
$ws.Range("L123").Value = 42743.332
$ws.Range("N123").Value = -52543.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 20007.5
$ws.Range("I21").Value = 10015
$ws.Range("K21").Value = 10015
$ws.Range("M21").Value = -9641
$ws.Range("H32").Value = 10627.61
$ws.Range("I32").Value = 7079.7256
$ws.Range("J32").Value = 17586.924
$ws.Range("K32").Value = 7079.7256
$ws.Range("L32").Value = 17586.924
$ws.Range("M32").Value = -6792.7256
$ws.Range("N32").Value = -18160.924
$ws.Range("H122").Value = 2543.5
$ws.Range("I122").Value = 1652.5333
$ws.Range("K122").Value = 4957.5999
$ws.Range("M122").Value = -2507.5999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1944.75
$ws.Range("I94").Value = 1661.2727
$ws.Range("J94").Value = 2984.1667
$ws.Range("K94").Value = 1661.2727
$ws.Range("L94").Value = 2984.1667
$ws.Range("M94").Value = -1210.2727
$ws.Range("N94").Value = -3886.1667
$ws.Range("H99").Value = 3076.36
$ws.Range("I99").Value = 1071.4615
$ws.Range("J99").Value = 5248.3335
$ws.Range("K99").Value = 1071.4615
$ws.Range("L99").Value = 5248.3335
$ws.Range("M99").Value = 426.5385000000001
$ws.Range("N99").Value = -8244.333500000001
$ws.Range("H107").Value = 1362.8125
$ws.Range("I107").Value = 1298.3846
$ws.Range("J107").Value = 1642
$ws.Range("K107").Value = 1298.3846
$ws.Range("L107").Value = 1642
$ws.Range("M107").Value = 621.6153999999999
$ws.Range("N107").Value = -5482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2689.4443
$ws.Range("I10").Value = 867.6667
$ws.Range("J10").Value = 6333
$ws.Range("K10").Value = 867.6667
$ws.Range("L10").Value = 6333
$ws.Range("M10").Value = -728.6667
$ws.Range("N10").Value = -6611
$ws.Range("H31").Value = 2668.1702
$ws.Range("I31").Value = 1178.5588
$ws.Range("J31").Value = 6564.077
$ws.Range("K31").Value = 1178.5588
$ws.Range("L31").Value = 6564.077
$ws.Range("M31").Value = -883.5588
$ws.Range("N31").Value = -7154.077
$ws.Range("H34").Value = 2668.1702
$ws.Range("I34").Value = 1178.5588
$ws.Range("J34").Value = 6564.077
$ws.Range("K34").Value = 1178.5588
$ws.Range("L34").Value = 6564.077
$ws.Range("M34").Value = -976.5588
$ws.Range("N34").Value = -6968.077
$ws.Range("H99").Value = 5801.5557
$ws.Range("I99").Value = 3080
$ws.Range("J99").Value = 9203.5
$ws.Range("K99").Value = 3080
$ws.Range("L99").Value = 9203.5
$ws.Range("M99").Value = -1582
$ws.Range("N99").Value = -12199.5
$ws.Range("H103").Value = 27503.125
$ws.Range("I103").Value = 10675
$ws.Range("J103").Value = 37600
$ws.Range("K103").Value = 10675
$ws.Range("L103").Value = 37600
$ws.Range("M103").Value = -9503
$ws.Range("N103").Value = -39944
$ws.Range("H107").Value = 1442
$ws.Range("I107").Value = 1442
$ws.Range("K107").Value = 1442
$ws.Range("M107").Value = 478
$ws.Range("H126").Value = 5801.5557
$ws.Range("I126").Value = 3080
$ws.Range("J126").Value = 9203.5
$ws.Range("K126").Value = 9240
$ws.Range("L126").Value = 27610.5
$ws.Range("M126").Value = -6770
$ws.Range("N126").Value = -32550.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 158.16667
$ws.Range("I18").Value = 158.16667
$ws.Range("K18").Value = 474.50001
$ws.Range("M18").Value = -305.50001
$ws.Range("H45").Value = 1136.1428
$ws.Range("J45").Value = 1142.1666
$ws.Range("L45").Value = 3426.4998
$ws.Range("N45").Value = -4490.4998
$ws.Range("H100").Value = 4323.75
$ws.Range("J100").Value = 4923.5713
$ws.Range("L100").Value = 14770.7139
$ws.Range("N100").Value = -16392.7139
$ws.Range("H103").Value = 1190.5555
$ws.Range("I103").Value = 785.8333
$ws.Range("J103").Value = 2000
$ws.Range("K103").Value = 2357.4999
$ws.Range("L103").Value = 6000
$ws.Range("M103").Value = -1478.4999
$ws.Range("N103").Value = -7758
$ws.Range("H105").Value = 4991.4287
$ws.Range("J105").Value = 4991.4287
$ws.Range("L105").Value = 14974.2861
$ws.Range("N105").Value = -20216.2861
$ws.Range("H108").Value = 2364.7144
$ws.Range("I108").Value = 2364.7144
$ws.Range("K108").Value = 7094.1432
$ws.Range("M108").Value = -4214.1432
$ws.Range("H112").Value = 4009.0715
$ws.Range("J112").Value = 4358.3335
$ws.Range("L112").Value = 13075.0005
$ws.Range("N112").Value = -15291.0005
$ws.Range("H117").Value = 5481.375
$ws.Range("I117").Value = 1430.25
$ws.Range("J117").Value = 9532.5
$ws.Range("K117").Value = 4290.75
$ws.Range("L117").Value = 28597.5
$ws.Range("M117").Value = -848.75
$ws.Range("N117").Value = -35481.5
$ws.Range("H121").Value = 1942.6923
$ws.Range("I121").Value = 304.83334
$ws.Range("J121").Value = 2109.2542
$ws.Range("K121").Value = 914.5000200000001
$ws.Range("L121").Value = 6327.7626
$ws.Range("M121").Value = 395.4999799999999
$ws.Range("N121").Value = -8947.7626
$ws.Range("H125").Value = 8450
$ws.Range("J125").Value = 9940
$ws.Range("L125").Value = 29820
$ws.Range("N125").Value = -39660
$ws.Range("H130").Value = 2857.375
$ws.Range("I130").Value = 1445
$ws.Range("K130").Value = 4335
$ws.Range("M130").Value = 685
$ws.Range("H138").Value = 3647.1428
$ws.Range("I138").Value = 3007.5
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 9022.5
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = -3882.5
$ws.Range("N138").Value = -23780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6428.4814
$ws.Range("I70").Value = 5773.45
$ws.Range("J70").Value = 8300
$ws.Range("K70").Value = 5773.45
$ws.Range("L70").Value = 8300
$ws.Range("M70").Value = -5503.45
$ws.Range("N70").Value = -8840
$ws.Range("H73").Value = 6428.4814
$ws.Range("I73").Value = 5773.45
$ws.Range("J73").Value = 8300
$ws.Range("K73").Value = 5773.45
$ws.Range("L73").Value = 8300
$ws.Range("M73").Value = -4837.45
$ws.Range("N73").Value = -10172
$ws.Range("H102").Value = 3033.9092
$ws.Range("I102").Value = 2418.3157
$ws.Range("K102").Value = 2418.3157
$ws.Range("M102").Value = -796.3157000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 10104066
$ws.Range("I93").Value = 13891465
$ws.Range("J93").Value = 4334.6665
$ws.Range("K93").Value = 13891465
$ws.Range("L93").Value = 4334.6665
$ws.Range("M93").Value = -13890217
$ws.Range("N93").Value = -6830.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 10000
$ws.Range("J12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("N12").Value = -10284
